$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 28 (pushes old rows 28..91 down to 29..92),
# matching the row above (27) exactly for formatting, then fill it in.
$ws.Rows.Item(28).Insert()

# Copy formatting (fill/border/font/alignment) from row 27, which is
# identical in style to what row 28 should look like (s=3 col A, s=4 rest).
$ws.Range("A27:J27").Copy()
$ws.Range("A28:J28").PasteSpecial(-4122)
$ws.Rows.Item(28).RowHeight = 75

# New row 28 content (Title / Step Action changed; the rest duplicates
# what used to be row 28, now shifted to row 29).
$ws.Range("A28").Value = "MatrixWeb: Navigation BurgerMenuNumbers_Verifica aggancio New Business Ultra Impresa"
$ws.Range("B28").Value = "Verifica aggancio New Business Ultra Impresa"
$ws.Range("C28").Value = "Si accede a Numbers, click burgerMenu e verifica atterraggio della pagina"
$ws.Range("D28").Value = "Pusateri Kevin (Leased Employed)"
$ws.Range("E28").Value = "Design"
$ws.Range("F28").Value = "Planned"
$ws.Range("G28").Value = "Sinistri"
$ws.Range("H28").Value = "Automation"
$ws.Range("I28").Value = "Test Factory"
$ws.Range("J28").Value = "Allianz Projects\Digital Interaction\Allianz Matrix Web\Numbers"

# Update view state to match the authored selection/scroll position.
$ws.Application.ActiveWindow.ScrollRow = 23
$ws.Range("J29").Select()
